# Reorders the "Med"/"Mean" derived-statistic columns of the corr3mvM
# correlation-matrix sheet so that each group reads
# mean, iqr, var, rmean, rvar (instead of iqr, mean, rmean, rvar, var).
#
# The sheet is a symmetric correlation matrix: row 1 (B1:U1) and column A
# (A2:A21) both carry the same ordered list of variable names, and the
# B2:U21 block holds the correlation of row-variable vs column-variable.
# Reordering the variable list therefore requires permuting BOTH the
# headers and the rows/columns of the data block together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old (current) and new (target) variable order.
$oldOrder = @(
    "sp500", "VIXCLS",
    "iqrMed", "meanMed", "rmeanMed", "rvarMed", "varMed",
    "iqrMean", "meanMean", "rmeanMean", "rvarMean", "varMean",
    "kurtEstMed", "meanEstMed", "skewEstMed", "varEstMed",
    "kurtEstMean", "meanEstMean", "skewEstMean", "varEstMean"
)

$newOrder = @(
    "sp500", "VIXCLS",
    "meanMed", "iqrMed", "varMed", "rmeanMed", "rvarMed",
    "meanMean", "iqrMean", "varMean", "rmeanMean", "rvarMean",
    "kurtEstMed", "meanEstMed", "skewEstMed", "varEstMed",
    "kurtEstMean", "meanEstMean", "skewEstMean", "varEstMean"
)

$n = $newOrder.Count

# For each position in the new order, find where that variable currently
# lives (0-based) so we know which old row/column to pull data from.
$oldIndexOf = @{}
for ($i = 0; $i -lt $n; $i++) {
    $oldIndexOf[$oldOrder[$i]] = $i
}
$srcPos = New-Object 'int[]' $n
for ($i = 0; $i -lt $n; $i++) {
    $srcPos[$i] = $oldIndexOf[$newOrder[$i]]
}

# Pull the current header row, row-label column and full data block in one
# shot each (columns B..U = variables 0..19, rows 2..21 = variables 0..19).
$headerRange = $ws.Range("B1:U1")
$labelRange  = $ws.Range("A2:A21")
$dataRange   = $ws.Range("B2:U21")

$headerVals = $headerRange.Value2
$labelVals  = $labelRange.Value2
$dataVals   = $dataRange.Value2

$newHeader = New-Object 'object[,]' 1, $n
$newLabel  = New-Object 'object[,]' $n, 1
$newData   = New-Object 'object[,]' $n, $n

for ($i = 0; $i -lt $n; $i++) {
    $si = $srcPos[$i] + 1   # 1-based source column/row index
    $newHeader[0, $i] = $headerVals[1, $si]
    $newLabel[$i, 0]  = $labelVals[$si, 1]
}

for ($r = 0; $r -lt $n; $r++) {
    $sr = $srcPos[$r] + 1
    for ($c = 0; $c -lt $n; $c++) {
        $sc = $srcPos[$c] + 1
        $newData[$r, $c] = $dataVals[$sr, $sc]
    }
}

$headerRange.Value2 = $newHeader
$labelRange.Value2  = $newLabel
$dataRange.Value2   = $newData
